# The workbook gained a new weekly price-report row. The new record is
# inserted as row 51 (a new "Arveja Verde" quote dated 2022-02-11, serial
# 44603), which pushes all the previously-existing rows 51-93 down by one
# (to 52-94). The sheet's used range therefore grows from A1:R93 to A1:R94.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 51, shifting rows 51:93 down to 52:94.
$ws.Rows(51).Insert()

# Populate the newly inserted row 51 with the new record's data.
$ws.Range("A51").Value = 10
$ws.Range("B51").Value = "Vega Modelo de Temuco"
$ws.Range("C51").Value = "La Araucanía"
$ws.Range("D51").Value = 44603
$ws.Range("E51").Value = 9
$ws.Range("F51").Value = 100112022
$ws.Range("G51").Value = "Arveja Verde"
$ws.Range("H51").Value = "Sin especificar"
$ws.Range("I51").Value = "Primera"
$ws.Range("J51").Value = 185
$ws.Range("K51").Value = 22000
$ws.Range("L51").Value = 22000
$ws.Range("M51").Value = 22000
$ws.Range("N51").Value = "$/malla 25 kilos"
$ws.Range("O51").Value = "Región de La Araucanía"
$ws.Range("P51").Value = 880
$ws.Range("Q51").Value = 25
$ws.Range("R51").Value = "Hortaliza"
